$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 already carries the right date/time styles (s=1 / s=3) -
# the commit only changes its date (a stray 2013-04-05 becomes 2013-11-05,
# i.e. serial 41369 -> 41583); the time value is unchanged.
$ws.Range("A24").Value = 41583
$ws.Range("B24").Value = 0.0625

# Two new rows are appended (25 and 26) continuing the time sheet.
# Duplicate row 24's formatting via Copy + Insert(shift-down) so the new
# rows land on the exact same cell styles as the rest of the column
# (rather than picking up a blank/default style), then fill in the
# real values for the new entries.
$ws.Range("A24:B24").Copy()
$ws.Range("A25:B25").Insert(-4121)
$ws.Range("A24:B24").Copy()
$ws.Range("A26:B26").Insert(-4121)

$ws.Range("A25").Value = 41584
$ws.Range("B25").Value = 0.090277777777777776

$ws.Range("A26").Value = 41588
$ws.Range("B26").Value = 0.097222222222222224

# Match the recorded selection after the edit.
# (cast to void - Range.Select returns a Boolean that would otherwise be
# echoed to the output stream)
[void]$ws.Range("C26").Select()

Write-Output "done"
